$d = $word.ActiveDocument

# 1. Relocate the "_GoBack" bookmark (Word automatically drops this bookmark at
#    the location of the most recent edit). Remove it from its old spot (end of
#    the "Christopher Mendez" byline paragraph) ...
$bookmarks = $d.Bookmarks
if ($bookmarks.Exists("_GoBack")) {
    $bookmarks.Item("_GoBack").Delete()
}

# 2. Insert the new sentence about additional data sources (Metacritic scores /
#    sentiment dictionary) right after the existing sentence about datasets.
$insertPoint = $d.Content
[void]$insertPoint.Find.Execute("training the model, and we kept 5 datasets in order to test the accuracy of the model. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($insertPoint.Find.Found) {
    $insertPoint.Collapse(0)
    $newSentence = "In addition to tweets, Metacritic scores and a sentiment dictionary were also collected as data. The Metacritic scores were collected by hand and entered in a text file to be used as input when building the model. In addition, the sentiment dictionary was imported from a text created by Finn Arup Nielsen which has a list of words and phrases given sentiment scores from 5 (being very positive) to -5 (being very negative)."
    $insertPoint.InsertAfter($newSentence)
    $insertPoint.Font.Spacing = 0.25
}

# 3. ... and drop "_GoBack" back at the point mid-sentence where the (quick) edit
#    was made, matching where Word would leave it after the author's last keystroke.
$goBackPoint = $d.Content
[void]$goBackPoint.Find.Execute("The final regression m", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($goBackPoint.Find.Found) {
    $goBackPoint.Collapse(0)
    $bookmarks.Add("_GoBack", $goBackPoint)
}
